$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "Planilha1" worksheet (data dictionary) between
#    "SaldoAnterior" and "Opções"
# ---------------------------------------------------------------------
$saldo = $wb.Worksheets.Item("SaldoAnterior")
$ws = $wb.Worksheets.Add($null, $saldo)
$ws.Name = "Planilha1"

# ---------------------------------------------------------------------
# 2. Populate the data dictionary content
# ---------------------------------------------------------------------

# Title
$ws.Range("B2").Value = "Tabela: SaldoAnterior"

# Column headers
$ws.Range("B4").Value = "Coluna"
$ws.Range("C4").Value = "Tipo de Dado"
$ws.Range("D4").Value = "Descrição"

# Row 5 - Banco_ID
$ws.Range("B5").Value = "Banco_ID"
$ws.Range("C5").Value = "Inteiro (int)"
$ws.Range("C5").Characters(10, 3).Font.Name = "Arial Unicode MS"
$ws.Range("C5").Characters(10, 3).Font.Size = 10
$ws.Range("D5").Value = "Identificador do banco, utilizado como chave estrangeira para relacionar à tabela Bancos."
$ws.Range("D5").Characters(70, 6).Font.Bold = $true

# Row 6 - Valor
$ws.Range("B6").Value = "Valor"
$ws.Range("C6").Value = "Número decimal (float)"
$ws.Range("C6").Characters(17, 5).Font.Name = "Arial Unicode MS"
$ws.Range("C6").Characters(17, 5).Font.Size = 10
$ws.Range("D6").Value = "Valor do saldo inicial do banco antes do início dos lançamentos de movimentos. Pode ser positivo ou negativo."

# ---------------------------------------------------------------------
# 3. Formatting
# ---------------------------------------------------------------------

# Column widths
$ws.Columns("B").ColumnWidth = 25.7109375
$ws.Columns("C").ColumnWidth = 25.7109375
$ws.Columns("D").ColumnWidth = 50.7109375

# Row heights for wrapped description rows
$ws.Rows("5").RowHeight = 30
$ws.Rows("6").RowHeight = 45

# Title style: centered, shaded fill
$ws.Range("B2").HorizontalAlignment = -4108
$ws.Range("B2").Interior.ThemeColor = 5

# Header row style: bold, centered, shaded fill, boxed border, wrap
$hdr = $ws.Range("B4:D4")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108
$hdr.WrapText = $true
$hdr.Interior.ThemeColor = 5
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

# Column-name cells (B5, B6): bold, centered, boxed border, wrap
$names = $ws.Range("B5:B6")
$names.Font.Bold = $true
$names.HorizontalAlignment = -4108
$names.VerticalAlignment = -4108
$names.WrapText = $true
$names.Borders.LineStyle = 1
$names.Borders.Weight = 2

# Data-type cells (C5, C6): centered, boxed border, wrap
$types = $ws.Range("C5:C6")
$types.HorizontalAlignment = -4108
$types.VerticalAlignment = -4108
$types.WrapText = $true
$types.Borders.LineStyle = 1
$types.Borders.Weight = 2

# Description cells (D5, D6): vertical centered, boxed border, wrap
$descr = $ws.Range("D5:D6")
$descr.VerticalAlignment = -4108
$descr.WrapText = $true
$descr.Borders.LineStyle = 1
$descr.Borders.Weight = 2

$ws.Range("A1").Select()

# ---------------------------------------------------------------------
# 4. View settings
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.DisplayGridlines = $false

# Restore SaldoAnterior as the active tab/selection (matches source file)
$saldo.Activate()
$saldo.Range("D16").Select()
